# Update countries & provincias Spain
# Refresh the "Pais" COVID stats sheet: new timestamp, updated case counts for
# several countries, and re-sorted rows where a country's total overtook its
# neighbour in the (descending) ranking so the two rows trade places.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 7 de Octubre de 2020 a las 22:38"

# --- Estados Unidos (row 4) --------------------------------------------
$ws.Range("B4").Value = 7757830
$ws.Range("C4").Value = 32101
$ws.Range("D4").Value = 4964455
$ws.Range("E4").Value = 2576951
$ws.Range("G4").Value = 572
$ws.Range("H4").Value = 216424

# --- India (row 5) ------------------------------------------------------
$ws.Range("B5").Value = 6832906
$ws.Range("C5").Value = 78727
$ws.Range("D5").Value = 5824415
$ws.Range("E5").Value = 902939
$ws.Range("G5").Value = 961
$ws.Range("H5").Value = 105552

# --- Alemania (row 26) ---------------------------------------------------
$ws.Range("B26").Value = 311061
$ws.Range("C26").Value = 3942
$ws.Range("E26").Value = 33709

# --- Israel (row 27) ------------------------------------------------------
$ws.Range("B27").Value = 281481
$ws.Range("C27").Value = 4455
$ws.Range("D27").Value = 216613
$ws.Range("E27").Value = 63044
$ws.Range("G27").Value = 27
$ws.Range("H27").Value = 1824

# --- Costa Rica (row 51) --------------------------------------------------
$ws.Range("B51").Value = 83497
$ws.Range("C51").Value = 1355
$ws.Range("D51").Value = 50295
$ws.Range("E51").Value = 32178
$ws.Range("G51").Value = 20
$ws.Range("H51").Value = 1024

# --- Cabo Verde (row 118) -------------------------------------------------
$ws.Range("B118").Value = 6624
$ws.Range("C118").Value = 106
$ws.Range("D118").Value = 5684
$ws.Range("E118").Value = 869
$ws.Range("G118").Value = 2
$ws.Range("H118").Value = 71

# --- Suazilandia (row 123) ------------------------------------------------
$ws.Range("B123").Value = 5617
$ws.Range("C123").Value = 19
$ws.Range("D123").Value = 5196
$ws.Range("E123").Value = 308
$ws.Range("G123").Value = 1
$ws.Range("H123").Value = 113

# --- Siria overtakes Sri Lanka: rows 135/136 swap ------------------------
$ws.Range("A135").Value = "Siria"
$ws.Range("B135").Value = 4504
$ws.Range("C135").Value = 47
$ws.Range("D135").Value = 1198
$ws.Range("E135").Value = 3094
$ws.Range("G135").Value = 3
$ws.Range("H135").Value = 212

$ws.Range("A136").Value = "Sri Lanka"
$ws.Range("B136").Value = 4459
$ws.Range("C136").Value = 207
$ws.Range("D136").Value = 3274
$ws.Range("E136").Value = 1172
$ws.Range("H136").Value = 13

# --- Republica del Chad (row 166) -----------------------------------------
$ws.Range("B166").Value = 1251
$ws.Range("C166").Value = 13
$ws.Range("D166").Value = 1090
$ws.Range("E166").Value = 72
$ws.Range("G166").Value = 1
$ws.Range("H166").Value = 89

# --- Curazao overtakes Comoras: rows 178/179 swap -------------------------
$ws.Range("A178").Value = "Curazao"
$ws.Range("B178").Value = 505
$ws.Range("C178").Value = 29
$ws.Range("D178").Value = 270
$ws.Range("E178").Value = 234
$ws.Range("H178").Value = 1

$ws.Range("A179").Value = "Comoras"
$ws.Range("B179").Value = 491
$ws.Range("D179").Value = 468
$ws.Range("E179").Value = 16
$ws.Range("H179").Value = 7

# --- Barbados (row 191) ----------------------------------------------------
$ws.Range("B191").Value = 203
$ws.Range("C191").Value = 3
$ws.Range("E191").Value = 14

# --- Santa Lucia overtakes Nueva Caledonia: rows 207/208 swap (stats tied) -
$ws.Range("A207").Value = "Santa Lucia"
$ws.Range("A208").Value = "Nueva Caledonia"

# --- Islas Malvinas overtakes Montserrat: rows 215/216 swap ----------------
$ws.Range("A215").Value = "Islas Malvinas"
$ws.Range("D215").Value = 13
$ws.Range("H215").Value = 0

$ws.Range("A216").Value = "Montserrat"
$ws.Range("D216").Value = 12
$ws.Range("H216").Value = 1
